{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, etc.)\n// in bold + corporate-blue color (#2C3E50) across the resume's achievement\n// and experience bullet points.\n//\n// Each entry identifies a target paragraph by a unique substring (so we\n// never touch the same numbers that also appear elsewhere in the document,\n// e.g. in the PROFESSIONAL SUMMARY or KEY PROJECTS sections), and lists --\n// in left-to-right order -- the metric substrings inside that paragraph\n// which must become bold + colored.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\nconst work = [\n  {\n    anchor: \"Discovered systematic race coding errors\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    anchor: \"Utilized advanced sampling methods\",\n    metrics: [\"\\u00B14.2%\", \"\\u00B12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    anchor: \"Trigonometric algorithm for boundary estimation\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    anchor: \"Built real-time FEC analysis systems\",\n    metrics: [\"$2\"],\n  },\n  {\n    anchor: \"Modernized legacy ETL processes\",\n    metrics: [\"57%\"],\n  },\n  {\n    anchor: \"Algorithmic innovation: Pioneered trigonometric\",\n    metrics: [\"73.5%\"],\n  },\n  {\n    anchor: \"$4.7M savings enabled nonprofit access\",\n    metrics: [\"$4.7M\"],\n  },\n  {\n    anchor: \"178% accuracy improvement\",\n    metrics: [\"178%\"],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const item of paragraphs.items) {\n  item.load(\"text\");\n}\nawait context.sync();\n\nfor (const job of work) {\n  const target = paragraphs.items.find((p) => p.text.indexOf(job.anchor) !== -1);\n\n  if (!target) {\n    console.log(\"Anchor paragraph not found: \" + job.anchor);\n    continue;\n  }\n\n  for (const metric of job.metrics) {\n    const results = target.search(metric, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length === 0) {\n      console.log(\"Metric not found in paragraph: \" + metric + \" / \" + job.anchor);\n      continue;\n    }\n\n    const hit = results.items[0];\n    hit.font.bold = true;\n    hit.font.color = HIGHLIGHT_COLOR;\n    await context.sync();\n  }\n}\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, etc.)\n# in bold + corporate-blue color across the resume's achievement and\n# experience bullet points.\n#\n# Each entry below identifies a target paragraph by a unique substring\n# (so we never touch the same numbers that also appear elsewhere in the\n# document, e.g. in the PROFESSIONAL SUMMARY or KEY PROJECTS sections),\n# and lists -- in left-to-right order -- the metric substrings inside\n# that paragraph which must become bold + colored 2C3E50.\n\n$d = $word.ActiveDocument\n\n# Bold + color (RGB 2C3E50, stored as BGR-ordered OLE color long)\n$highlightColor = 0x50 * 65536 + 0x3E * 256 + 0x2C\n\n$work = @(\n    @{ Anchor = \"Discovered systematic race coding errors\"; Metrics = @(\"23%\", \"64%\") },\n    @{ Anchor = \"Utilized advanced sampling methods\"; Metrics = @([char]0x00B1 + \"4.2%\", [char]0x00B1 + \"2.1%\", \"71%\", \"87%\") },\n    @{ Anchor = \"Trigonometric algorithm for boundary estimation\"; Metrics = @(\"73.5%\", \"`$4.7M\") },\n    @{ Anchor = \"Built real-time FEC analysis systems\"; Metrics = @(\"`$2\") },\n    @{ Anchor = \"Modernized legacy ETL processes\"; Metrics = @(\"57%\") },\n    @{ Anchor = \"Algorithmic innovation: Pioneered trigonometric\"; Metrics = @(\"73.5%\") },\n    @{ Anchor = \"`$4.7M savings enabled nonprofit access\"; Metrics = @(\"`$4.7M\") },\n    @{ Anchor = \"178% accuracy improvement\"; Metrics = @(\"178%\") }\n)\n\nforeach ($item in $work) {\n    $anchor = $item.Anchor\n    $targetIndex = -1\n\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $pText = $d.Paragraphs.Item($i).Range.Text\n        if ($pText -like (\"*\" + $anchor + \"*\")) {\n            $targetIndex = $i\n            break\n        }\n    }\n\n    if ($targetIndex -eq -1) {\n        Write-Output (\"Anchor paragraph not found: \" + $anchor)\n        continue\n    }\n\n    foreach ($metric in $item.Metrics) {\n        $para = $d.Paragraphs.Item($targetIndex)\n        $rng = $para.Range\n        $find = $rng.Find\n        $find.ClearFormatting()\n        $find.Text = $metric\n        $find.MatchCase = $true\n        $find.MatchWildcards = $false\n        [void]$find.Execute()\n\n        if ($find.Found) {\n            $rng.Font.Bold = $true\n            $rng.Font.Color = $highlightColor\n        } else {\n            Write-Output (\"Metric not found in paragraph: \" + $metric + \" / \" + $anchor)\n        }\n    }\n}\n"}
